$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.269585132598877
$ws.Range("B1").Value = 2.512616872787476
$ws.Range("C1").Value = 5.067704677581787
$ws.Range("D1").Value = 2.86467432975769
$ws.Range("E1").Value = 1.096312046051025
